$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old data block first so stale cells do not linger ---
$ws.Range("A1:H5").ClearContents()

# --- Row 1: header ---
$ws.Range("A1").Value = "id"

# --- Apply the "style 1" (vertical alignment reset from the workbook's
#     default vertical=center) to every cell that carries it in the target ---
$ws.Range("A2:C2").VerticalAlignment = -4107
$ws.Range("A3:H3").VerticalAlignment = -4107
$ws.Range("A4:D4").VerticalAlignment = -4107

# --- Helper: write a value into a cell as TEXT (no numeric auto-coercion)
#     by staging it in a scratch cell formatted as Text, then pasting
#     values-only into the destination (keeps destination's own style). ---
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-TextValue($rangeAddr, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)
}

# Row 2
Set-TextValue "A2" "HMDB0001311"
Set-TextValue "B2" "1.4"
Set-TextValue "C2" "1.41"

# Row 3
Set-TextValue "A3" "HMDB0001485"
Set-TextValue "B3" "0.88"
Set-TextValue "C3" "0.89"
Set-TextValue "D3" "0.91"
Set-TextValue "E3" "1.27"
Set-TextValue "F3" "1.3"
Set-TextValue "G3" "1.31"
Set-TextValue "H3" "1.33"

# Row 4
Set-TextValue "A4" "HMDB0001494"
Set-TextValue "B4" "1.91"
Set-TextValue "C4" "2.1"
Set-TextValue "D4" "2.1"
$ws.Range("E4").Value = 8.1999999999999993

# Row 5 (plain/default style)
$ws.Range("A5").Value = "HMDB0002361"
$ws.Range("B5").Value = 1.28

# --- Clean up the scratch cell/column entirely ---
$scratch.EntireColumn.Delete()

# --- Column width + selection + dimension cosmetics ---
# (ColumnWidth COM units are offset from the raw XML "width" by 5/7 = 0.714285...)
$ws.Columns("A").ColumnWidth = 15.285714285714286
$ws.Range("F17").Select()
